$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '51.879.26'
Set-TextValue 'E2' '  +0.07%  '
Set-TextValue 'D3' '2.778.03'
Set-TextValue 'E3' '  -2.10%  '
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '356.24'
Set-TextValue 'E5' '  +1.05%  '
Set-TextValue 'D6' '109.36'
Set-TextValue 'E6' '  -4.04%  '
Set-TextValue 'D7' '0.565'
Set-TextValue 'E7' '  +0.92%  '
Set-TextValue 'E8' '  +0.10%  '
Set-TextValue 'D9' '0.596'
Set-TextValue 'E9' '  -0.81%  '
Set-TextValue 'D10' '39.95'
Set-TextValue 'E10' '  -4.22%  '
Set-TextValue 'E11' '  -0.21%  '
Set-TextValue 'E12' '  +0.77%  '
Set-TextValue 'D13' '19.43'
Set-TextValue 'E13' '  -2.82%  '
Set-TextValue 'E14' '  -1.73%  '
Set-TextValue 'D15' '3.214.48'
Set-TextValue 'E15' '  -1.87%  '
Set-TextValue 'D16' '2.783.00'
Set-TextValue 'E16' '  -1.59%  '
Set-TextValue 'D17' '0.934'
Set-TextValue 'E17' '  +3.97%  '
Set-TextValue 'D18' '51.760.57'
Set-TextValue 'E18' '  +0.13%  '
Set-TextValue 'E19' '  +0.44%  '
Set-TextValue 'D20' '3.13'
Set-TextValue 'E20' '  -0.42%  '
Set-TextValue 'D21' '13.01'
Set-TextValue 'E21' '  -3.60%  '
Set-TextValue 'E22' '  -1.80%  '
Set-TextValue 'D23' '274.24'
Set-TextValue 'E23' '  +1.20%  '
Set-TextValue 'D24' '69.98'
Set-TextValue 'E24' '  +0.32%  '
Set-TextValue 'E25' '  -1.58%  '
Set-TextValue 'D26' '26.62'
Set-TextValue 'E26' '  -0.47%  '
Set-TextValue 'E27' '  -0.08%  '
Set-TextValue 'D28' '10.15'
Set-TextValue 'E28' '  -1.50%  '
Set-TextValue 'E29' '  +4.07%  '
Set-TextValue 'D30' '2.22'
Set-TextValue 'E30' '  -1.50%  '
Set-TextValue 'B31' 'VeChain'
Set-TextValue 'C31' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D31' '0.0466'
Set-TextValue 'E31' '  +3.62%  '
Set-TextValue 'B32' 'OKB'
Set-TextValue 'C32' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D32' '51.58'
Set-TextValue 'E32' '  +1.71%  '
Set-TextValue 'D33' '34.01'
Set-TextValue 'E33' '  +0.34%  '
Set-TextValue 'E34' '  -2.01%  '
Set-TextValue 'D35' '0.0845'
Set-TextValue 'E35' '  +2.08%  '
Set-TextValue 'D36' '5.25'
Set-TextValue 'E36' '  +7.31%  '
Set-TextValue 'E37' '  +0.02%  '
Set-TextValue 'E38' '  -0.19%  '
Set-TextValue 'D39' '18.11'
Set-TextValue 'E39' '  +0.43%  '
Set-TextValue 'D40' '2.00'
Set-TextValue 'E40' '  -4.02%  '
Set-TextValue 'E41' '  -0.33%  '
Set-TextValue 'D42' '2.52'
Set-TextValue 'E42' '  -1.61%  '
Set-TextValue 'B43' 'Monero'
Set-TextValue 'C43' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D43' '121.72'
Set-TextValue 'E43' '  -3.10%  '
Set-TextValue 'B44' 'WEMIXToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D44' '2.24'
Set-TextValue 'E44' '  -2.71%  '
Set-TextValue 'D45' '22.03'
Set-TextValue 'E45' '  -6.35%  '
Set-TextValue 'D46' '2.065.41'
Set-TextValue 'E46' '  -0.83%  '
Set-TextValue 'D47' '3.24'
Set-TextValue 'E47' '  -3.26%  '
Set-TextValue 'D48' '2.18'
Set-TextValue 'E48' '  -5.83%  '
Set-TextValue 'E49' '  -0.28%  '
Set-TextValue 'E50' '  -0.59%  '
Set-TextValue 'D51' '8.94'
Set-TextValue 'E51' '  +0.10%  '
